# Add data for 2021-10-04: update the "through 09-25" snapshot to "through 09-26"
# by refreshing the worksheet/tab name, the September row label, and the
# September + Total rows with the new counts/rates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / tab to reflect the new "as-of" date.
$ws.Name = "Through 2021-09-26"

# Update the row-11 ("September") label.
$ws.Range("A11").Value = "September (through 09-26)"

# Row 11 - September monthly figures (arrest_made, no_arrest_made, arrest_rate per year).
$ws.Range("C11").Value = 26
$ws.Range("D11").Value = 0.037
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 36
$ws.Range("G11").Value = 0.0769
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 62
$ws.Range("J11").Value = 0.0746
$ws.Range("L11").Value = 46
$ws.Range("M11").Value = 0.08
$ws.Range("O11").Value = 57
$ws.Range("P11").Value = 0.0952
$ws.Range("R11").Value = 96
$ws.Range("S11").Value = 0.0303
$ws.Range("U11").Value = 162
$ws.Range("V11").Value = 0.0122

# Row 12 - Total figures (arrest_made, no_arrest_made, arrest_rate per year).
$ws.Range("C12").Value = 191
$ws.Range("D12").Value = 0.1357
$ws.Range("E12").Value = 44
$ws.Range("F12").Value = 376
$ws.Range("G12").Value = 0.1048
$ws.Range("H12").Value = 50
$ws.Range("I12").Value = 568
$ws.Range("J12").Value = 0.0809
$ws.Range("L12").Value = 479
$ws.Range("M12").Value = 0.113
$ws.Range("O12").Value = 370
$ws.Range("P12").Value = 0.1019
$ws.Range("R12").Value = 833
$ws.Range("S12").Value = 0.0566
$ws.Range("U12").Value = 1158
$ws.Range("V12").Value = 0.0616
